$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 17.32696144875414
$ws.Range("P2").Value = 33.49229712085892
$ws.Range("Q2").Value = 34.39035951189135
$ws.Range("R2").Value = 35.41391701186709
$ws.Range("S2").Value = 84.92854605168267
$ws.Range("T2").Value = 85.02983824704027
$ws.Range("U2").Value = 85.64287068350005
$ws.Range("V2").Value = 86.88053710548009
$ws.Range("W2").Value = 88.07208021552728
$ws.Range("X2").Value = 81.86940173570979
$ws.Range("Y2").Value = 82.13676894029437
$ws.Range("Z2").Value = 80.89695578181423
$ws.Range("AA2").Value = 79.55134706089075
$ws.Range("AB2").Value = 80.93117727286746
$ws.Range("AC2").Value = 28.82934562369141
$ws.Range("AD2").Value = 32.98787764209279
$ws.Range("AE2").Value = 37.00066340813115
$ws.Range("AF2").Value = 36.16370408496176
$ws.Range("AG2").Value = 34.2327719916449
$ws.Range("AH2").Value = 26.15438006310633
$ws.Range("AI2").Value = 24.50581304935785
$ws.Range("AJ2").Value = 22.49461324689247
$ws.Range("AK2").Value = 20.92555656828924
$ws.Range("AL2").Value = 24.29788302425395
$ws.Range("AM2").Value = 32.81906623717911
$ws.Range("AN2").Value = 42.32333432585737
$ws.Range("AO2").Value = 55.96193739139504
$ws.Range("AP2").Value = 73.56913035839048
$ws.Range("AQ2").Value = 82.52753925372205
$ws.Range("AR2").Value = 83.26134351992908
$ws.Range("AS2").Value = 83.01595631059399
$ws.Range("AT2").Value = 81.83999657041021
$ws.Range("AU2").Value = 79.85621727345422
$ws.Range("AV2").Value = 79.02554276929158
$ws.Range("AW2").Value = 77.57320485635391
$ws.Range("AX2").Value = 81.35415645165392
$ws.Range("AY2").Value = 69.40637735576712
$ws.Range("AZ2").Value = 53.00663904691545
$ws.Range("BA2").Value = 45.77513722326848
$ws.Range("BB2").Value = 45.10795235763289
$ws.Range("BC2").Value = 43.83255968789291
$ws.Range("BD2").Value = 40.50003902636562
$ws.Range("BE2").Value = 36.5218648367924
$ws.Range("BF2").Value = 33.9639007549929
$ws.Range("BG2").Value = 32.27302526741257
$ws.Range("BH2").Value = 33.92527604494219
$ws.Range("BI2").Value = 39.08266757431116
$ws.Range("BJ2").Value = 47.74202283965818
$ws.Range("BK2").Value = 58.79598381103528
$ws.Range("BL2").Value = 70.3260918591609
$ws.Range("BM2").Value = 81.51640448691518
$ws.Range("BN2").Value = 88.13900052225496
$ws.Range("BO2").Value = 89.99783170844684
$ws.Range("BP2").Value = 89.63040455173902
$ws.Range("BQ2").Value = 87.99052447736831
$ws.Range("BR2").Value = 86.6643953354595
$ws.Range("BS2").Value = 85.34626569211129
$ws.Range("BT2").Value = 83.58954003208684
$ws.Range("BU2").Value = 79.61152998233877
$ws.Range("BV2").Value = 70.65853611356371
$ws.Range("BW2").Value = 60.96975810921349
$ws.Range("BX2").Value = 54.99345585638441
$ws.Range("BY2").Value = 52.47408503270037
$ws.Range("BZ2").Value = 50.50639557378098
$ws.Range("CA2").Value = 47.30395967948595
$ws.Range("CB2").Value = 44.74807277460272
$ws.Range("CC2").Value = 43.36823255407359
$ws.Range("CD2").Value = 44.00544538680705
$ws.Range("CE2").Value = 47.42308403909227
$ws.Range("CF2").Value = 53.21527461957481
$ws.Range("CG2").Value = 62.13974721639234
$ws.Range("CH2").Value = 72.36356628358456
$ws.Range("CI2").Value = 82.3227609323248
$ws.Range("CJ2").Value = 90.38777896610641
$ws.Range("CK2").Value = 94.33329832923295
$ws.Range("CL2").Value = 95.96390303607997
$ws.Range("CM2").Value = 95.58169013911288
$ws.Range("CN2").Value = 94.29985073763316
$ws.Range("CO2").Value = 92.65344011572105
$ws.Range("O3").Value = 5.784025224107451
$ws.Range("P3").Value = 11.50108520489411
$ws.Range("Q3").Value = 11.76010928367492
$ws.Range("R3").Value = 12.04415820648252
$ws.Range("S3").Value = 29.59292952729937
$ws.Range("T3").Value = 29.52589716769242
$ws.Range("U3").Value = 29.61009968350706
$ws.Range("V3").Value = 29.89509454435627
$ws.Range("W3").Value = 30.17410459035305
$ws.Range("X3").Value = 27.84526688590122
$ws.Range("Y3").Value = 27.85566953396089
$ws.Range("Z3").Value = 27.37743990476528
$ws.Range("AA3").Value = 26.89049819657439
$ws.Range("AB3").Value = 27.34357583690048
$ws.Range("AC3").Value = 8.705954996976994
$ws.Range("AD3").Value = 10.10701432457626
$ws.Range("AE3").Value = 11.46029034356551
$ws.Range("AF3").Value = 11.20449448888429
$ws.Range("AG3").Value = 10.58534913740484
$ws.Range("AH3").Value = 7.736786132177394
$ws.Range("AI3").Value = 7.211132554220864
$ws.Range("AJ3").Value = 6.565303171013909
$ws.Range("AK3").Value = 6.066311151160124
$ws.Range("AL3").Value = 7.206432756871537
$ws.Range("AM3").Value = 10.05470646830594
$ws.Range("AN3").Value = 13.22902501660861
$ws.Range("AO3").Value = 17.77432033424235
$ws.Range("AP3").Value = 23.63533570582094
$ws.Range("AQ3").Value = 26.62703756873486
$ws.Range("AR3").Value = 26.88934781505467
$ws.Range("AS3").Value = 26.82531197379906
$ws.Range("AT3").Value = 26.45073529009458
$ws.Range("AU3").Value = 25.80609568568834
$ws.Range("AV3").Value = 25.55176933448388
$ws.Range("AW3").Value = 25.08093824846767
$ws.Range("AX3").Value = 26.56134658933731
$ws.Range("AY3").Value = 22.60282563560431
$ws.Range("AZ3").Value = 17.16630094354478
$ws.Range("BA3").Value = 14.76946031358933
$ws.Range("BB3").Value = 14.52084483606968
$ws.Range("BC3").Value = 14.09760883111113
$ws.Range("BD3").Value = 12.99227085608661
$ws.Range("BE3").Value = 11.67327754544221
$ws.Range("BF3").Value = 10.82601150470918
$ws.Range("BG3").Value = 10.23951967673031
$ws.Range("BH3").Value = 10.78957552366638
$ws.Range("BI3").Value = 12.50288659125467
$ws.Range("BJ3").Value = 15.37823915528196
$ws.Range("BK3").Value = 19.04815951860363
$ws.Range("BL3").Value = 22.85527902510412
$ws.Range("BM3").Value = 26.56991088567782
$ws.Range("BN3").Value = 28.76936423178495
$ws.Range("BO3").Value = 29.38844178112322
$ws.Range("BP3").Value = 29.26857828384788
$ws.Range("BQ3").Value = 28.71099762107588
$ws.Range("BR3").Value = 28.27079161723522
$ws.Range("BS3").Value = 27.83234311625544
$ws.Range("BT3").Value = 27.24748711563495
$ws.Range("BU3").Value = 25.9249643155673
$ws.Range("BV3").Value = 22.93694240925779
$ws.Range("BW3").Value = 19.71781552692856
$ws.Range("BX3").Value = 17.72920866996908
$ws.Range("BY3").Value = 16.88650589141938
$ws.Range("BZ3").Value = 16.22632083986759
$ws.Range("CA3").Value = 15.14109832115481
$ws.Range("CB3").Value = 14.28472974384531
$ws.Range("CC3").Value = 13.8185379407105
$ws.Range("CD3").Value = 14.02161149114395
$ws.Range("CE3").Value = 15.14725283943443
$ws.Range("CF3").Value = 17.04369462712635
$ws.Range("CG3").Value = 19.99582171213567
$ws.Range("CH3").Value = 23.37936485788511
$ws.Range("CI3").Value = 26.67554621939716
$ws.Range("CJ3").Value = 29.3438071349347
$ws.Range("CK3").Value = 30.62283071322303
$ws.Range("CL3").Value = 31.15645123783662
$ws.Range("CM3").Value = 31.02241065263608
$ws.Range("CN3").Value = 30.58977726509923
$ws.Range("CO3").Value = 30.03588680263502
$ws.Range("O4").Value = 205061.1594476519
$ws.Range("P4").Value = 222142.4730103448
$ws.Range("Q4").Value = 230882.5684367957
$ws.Range("R4").Value = 239454.1523822866
$ws.Range("S4").Value = 271681.2850702453
$ws.Range("T4").Value = 280304.6673240738
$ws.Range("U4").Value = 288272.534184465
$ws.Range("V4").Value = 295482.331369767
$ws.Range("W4").Value = 302502.288260672
$ws.Range("X4").Value = 309654.2370926514
$ws.Range("Y4").Value = 328504.5937478673
$ws.Range("Z4").Value = 362914.0176681723
$ws.Range("AA4").Value = 401434.7498377926
$ws.Range("AB4").Value = 426956.0052739589
$ws.Range("AC4").Value = 412623.2821405843
$ws.Range("AD4").Value = 419473.0404047451
$ws.Range("AE4").Value = 425803.1874671596
$ws.Range("AF4").Value = 433316.7463622336
$ws.Range("AG4").Value = 442194.4608687173
$ws.Range("AH4").Value = 447940.5906740123
$ws.Range("AI4").Value = 459201.5110182888
$ws.Range("AJ4").Value = 471030.4852458419
$ws.Range("AK4").Value = 483363.883210192
$ws.Range("AL4").Value = 496260.1017703174
$ws.Range("AM4").Value = 510105.2979685136
$ws.Range("AN4").Value = 524100.3538404885
$ws.Range("AO4").Value = 538569.5468106343
$ws.Range("AP4").Value = 553398.6307598195
$ws.Range("AQ4").Value = 566927.4604220003
$ws.Range("AR4").Value = 579942.0568771884
$ws.Range("AS4").Value = 592404.3241843701
$ws.Range("AT4").Value = 604254.6440104519
$ws.Range("AU4").Value = 615460.8275696888
$ws.Range("AV4").Value = 628808.0773624503
$ws.Range("AW4").Value = 640932.6402820915
$ws.Range("AX4").Value = 652391.415992927
$ws.Range("AY4").Value = 661402.1107656367
$ws.Range("AZ4").Value = 669488.4616051657
$ws.Range("BA4").Value = 678592.4550751266
$ws.Range("BB4").Value = 689631.304189756
$ws.Range("BC4").Value = 701234.4146979725
$ws.Range("BD4").Value = 712873.6056027861
$ws.Range("BE4").Value = 724738.1012371634
$ws.Range("BF4").Value = 737199.6869281329
$ws.Range("BG4").Value = 751066.3084124668
$ws.Range("BH4").Value = 766060.2005501325
$ws.Range("BI4").Value = 781689.361074789
$ws.Range("BJ4").Value = 797888.4565502958
$ws.Range("BK4").Value = 814545.5786151728
$ws.Range("BL4").Value = 832183.1157554951
$ws.Range("BM4").Value = 850127.6115915377
$ws.Range("BN4").Value = 867509.1875116909
$ws.Range("BO4").Value = 884240.4455797914
$ws.Range("BP4").Value = 900585.016336814
$ws.Range("BQ4").Value = 917687.9293073762
$ws.Range("BR4").Value = 934834.0535152019
$ws.Range("BS4").Value = 951759.5536544924
$ws.Range("BT4").Value = 968377.6233152276
$ws.Range("BU4").Value = 984448.9326784801
$ws.Range("BV4").Value = 1000730.582157382
$ws.Range("BW4").Value = 1016930.275515371
$ws.Range("BX4").Value = 1033531.239642103
$ws.Range("BY4").Value = 1050606.718255855
$ws.Range("BZ4").Value = 1067856.481627514
$ws.Range("CA4").Value = 1086178.26938479
$ws.Range("CB4").Value = 1104934.146111237
$ws.Range("CC4").Value = 1124039.892779936
$ws.Range("CD4").Value = 1143620.279578571
$ws.Range("CE4").Value = 1163747.515935079
$ws.Range("CF4").Value = 1185620.601535605
$ws.Range("CG4").Value = 1208289.058001181
$ws.Range("CH4").Value = 1231337.15361812
$ws.Range("CI4").Value = 1254582.15345857
$ws.Range("CJ4").Value = 1277796.379358334
$ws.Range("CK4").Value = 1301742.348075641
$ws.Range("CL4").Value = 1325807.349000925
$ws.Range("CM4").Value = 1349671.426940263
$ws.Range("CN4").Value = 1373416.012986802
$ws.Range("CO4").Value = 1397048.622448626
$ws.Range("O5").Value = 1067.377959174623
$ws.Range("P5").Value = 6648.888516293018
$ws.Range("Q5").Value = 6692.246893491353
$ws.Range("R5").Value = 6745.721410697191
$ws.Range("S5").Value = 21100.68122010895
$ws.Range("T5").Value = 21176.06505995277
$ws.Range("U5").Value = 21262.13566865758
$ws.Range("V5").Value = 21357.85007919374
$ws.Range("W5").Value = 21461.49684719602
$ws.Range("X5").Value = 20629.85802051304
$ws.Range("Y5").Value = 20741.94791610423
$ws.Range("Z5").Value = 20854.91918232272
$ws.Range("AA5").Value = 20970.46251387635
$ws.Range("AB5").Value = 21093.2834463315
$ws.Range("AC5").Value = 4537.049306351149
$ws.Range("AD5").Value = 4698.788133019242
$ws.Range("AE5").Value = 4902.551813119204
$ws.Range("AF5").Value = 5170.879167893521
$ws.Range("AG5").Value = 5530.336119257551
$ws.Range("AH5").Value = 2851.342796074126
$ws.Range("AI5").Value = 3475.238633803094
$ws.Range("AJ5").Value = 4264.744866426124
$ws.Range("AK5").Value = 5230.707368983592
$ws.Range("AL5").Value = 6371.318993579087
$ws.Range("AM5").Value = 7669.828574119255
$ws.Range("AN5").Value = 9093.287760228441
$ws.Range("AO5").Value = 10592.49592975919
$ws.Range("AP5").Value = 12103.38360334139
$ws.Range("AQ5").Value = 13550.25602808671
$ws.Range("AR5").Value = 14851.26526547222
$ws.Range("AS5").Value = 15925.83542609151
$ws.Range("AT5").Value = 16703.37657156867
$ws.Range("AU5").Value = 17131.99543433898
$ws.Range("AV5").Value = 19315.12397257435
$ws.Range("AW5").Value = 19599.50376163373
$ws.Range("AX5").Value = 18942.58379150126
$ws.Range("AY5").Value = 18005.58718540855
$ws.Range("AZ5").Value = 16866.60398859566
$ws.Range("BA5").Value = 15615.46526487737
$ws.Range("BB5").Value = 14036.72143701217
$ws.Range("BC5").Value = 12840.47194213542
$ws.Range("BD5").Value = 11799.3662097145
$ws.Range("BE5").Value = 10980.28554407684
$ws.Range("BF5").Value = 10429.53288091462
$ws.Range("BG5").Value = 9863.984051103012
$ws.Range("BH5").Value = 9891.492005218415
$ws.Range("BI5").Value = 10182.46230949343
$ws.Range("BJ5").Value = 10696.74899011261
$ws.Range("BK5").Value = 11385.39760665719
$ws.Range("BL5").Value = 11967.15951455483
$ws.Range("BM5").Value = 12850.60543466904
$ws.Range("BN5").Value = 13758.70705607637
$ws.Range("BO5").Value = 14647.71558781229
$ws.Range("BP5").Value = 15478.07476426947
$ws.Range("BQ5").Value = 16055.04005707327
$ws.Range("BR5").Value = 16669.95268777067
$ws.Range("BS5").Value = 17140.60663857036
$ws.Range("BT5").Value = 17452.37368185749
$ws.Range("BU5").Value = 17598.73985903091
$ws.Range("BV5").Value = 17425.61673350647
$ws.Range("BW5").Value = 17256.59451023968
$ws.Range("BX5").Value = 16955.60100273897
$ws.Range("BY5").Value = 16551.16944823026
$ws.Range("BZ5").Value = 16078.76750949685
$ws.Range("CA5").Value = 15414.40813912413
$ws.Range("CB5").Value = 14927.46346425973
$ws.Range("CC5").Value = 14492.97262965029
$ws.Range("CD5").Value = 14144.18859412355
$ws.Range("CE5").Value = 13905.83076709068
$ws.Range("CF5").Value = 13607.56061940195
$ws.Range("CG5").Value = 13623.49316991063
$ws.Range("CH5").Value = 13763.23287013216
$ws.Range("CI5").Value = 14013.20873787401
$ws.Range("CJ5").Value = 14353.60383527226
$ws.Range("CK5").Value = 14507.58598014426
$ws.Range("CL5").Value = 14954.43755713589
$ws.Range("CM5").Value = 15414.67122657487
$ws.Range("CN5").Value = 15862.83307569016
$ws.Range("CO5").Value = 16276.09460719779
